$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5525.2856
$ws.Range("J51").Value = 5340.5
$ws.Range("L51").Value = 5340.5
$ws.Range("N51").Value = -6308.5
$ws.Range("H58").Value = 2294
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 45000
$ws.Range("N58").Value = -45300
$ws.Range("H70").Value = 787163
$ws.Range("I70").Value = 2042078
$ws.Range("J70").Value = 2841.125
$ws.Range("K70").Value = 6126234
$ws.Range("L70").Value = 8523.375
$ws.Range("M70").Value = -6125964
$ws.Range("N70").Value = -9063.375
$ws.Range("H73").Value = 787163
$ws.Range("I73").Value = 2042078
$ws.Range("J73").Value = 2841.125
$ws.Range("K73").Value = 6126234
$ws.Range("L73").Value = 8523.375
$ws.Range("M73").Value = -6125298
$ws.Range("N73").Value = -10395.375
$ws.Range("H80").Value = 762454.0600000001
$ws.Range("J80").Value = 1543.8
$ws.Range("L80").Value = 4631.4
$ws.Range("N80").Value = -6627.4
$ws.Range("H83").Value = 762454.0600000001
$ws.Range("J83").Value = 1543.8
$ws.Range("L83").Value = 13894.2
$ws.Range("N83").Value = -23878.2
$ws.Range("H111").Value = 26112.555
$ws.Range("I111").Value = 1671
$ws.Range("J111").Value = 56664.5
$ws.Range("K111").Value = 5013
$ws.Range("L111").Value = 169993.5
$ws.Range("M111").Value = -1946
$ws.Range("N111").Value = -176127.5
$ws.Range("H116").Value = 83496670
$ws.Range("I116").Value = 50222004
$ws.Range("K116").Value = 50222004
$ws.Range("M116").Value = -50218562

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2859.5625
$ws.Range("I32").Value = 1521.0182
$ws.Range("K32").Value = 1521.0182
$ws.Range("M32").Value = -1234.0182
$ws.Range("H45").Value = 2029.5
$ws.Range("I45").Value = 1810.5555
$ws.Range("K45").Value = 1810.5555
$ws.Range("M45").Value = -1433.5555
$ws.Range("H80").Value = 48180.91
$ws.Range("J80").Value = 50499
$ws.Range("L80").Value = 50499
$ws.Range("N80").Value = -52495
$ws.Range("H83").Value = 48180.91
$ws.Range("J83").Value = 50499
$ws.Range("L83").Value = 151497
$ws.Range("N83").Value = -161481
$ws.Range("H101").Value = 53767
$ws.Range("J101").Value = 53767
$ws.Range("L101").Value = 53767
$ws.Range("N101").Value = -60257
$ws.Range("H110").Value = 62565080
$ws.Range("I110").Value = 76964180
$ws.Range("J110").Value = 169016.33
$ws.Range("K110").Value = 76964180
$ws.Range("L110").Value = 169016.33
$ws.Range("M110").Value = -76962135
$ws.Range("N110").Value = -173106.33
$ws.Range("H132").Value = 27779438
$ws.Range("I132").Value = 31251508
$ws.Range("K132").Value = 93754524
$ws.Range("M132").Value = -93751994

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 783.5
$ws.Range("I20").Value = 800
$ws.Range("J20").Value = 767
$ws.Range("K20").Value = 800
$ws.Range("L20").Value = 767
$ws.Range("M20").Value = -553
$ws.Range("N20").Value = -1261
$ws.Range("H86").Value = 7600
$ws.Range("I86").Value = 8533.333000000001
$ws.Range("K86").Value = 8533.333000000001
$ws.Range("M86").Value = -7410.333000000001
$ws.Range("H89").Value = 7600
$ws.Range("I89").Value = 8533.333000000001
$ws.Range("K89").Value = 42666.665
$ws.Range("M89").Value = -37050.665
$ws.Range("H94").Value = 4057.7273
$ws.Range("I94").Value = 4580.25
$ws.Range("J94").Value = 2664.3333
$ws.Range("K94").Value = 4580.25
$ws.Range("L94").Value = 2664.3333
$ws.Range("M94").Value = -4129.25
$ws.Range("N94").Value = -3566.3333
$ws.Range("H105").Value = 1899.3572
$ws.Range("I105").Value = 1771.909
$ws.Range("K105").Value = 1771.909
$ws.Range("M105").Value = -24.90900000000011
$ws.Range("H107").Value = 83541660
$ws.Range("I107").Value = 250000
$ws.Range("K107").Value = 250000
$ws.Range("M107").Value = -248080
$ws.Range("H134").Value = 3243.8823
$ws.Range("I134").Value = 3196.6562
$ws.Range("K134").Value = 9589.9686
$ws.Range("M134").Value = -7054.9686

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3027.1177
$ws.Range("I31").Value = 1482.0588
$ws.Range("J31").Value = 3542.1372
$ws.Range("K31").Value = 1482.0588
$ws.Range("L31").Value = 3542.1372
$ws.Range("M31").Value = -1187.0588
$ws.Range("N31").Value = -4132.1372
$ws.Range("H34").Value = 3027.1177
$ws.Range("I34").Value = 1482.0588
$ws.Range("J34").Value = 3542.1372
$ws.Range("K34").Value = 1482.0588
$ws.Range("L34").Value = 3542.1372
$ws.Range("M34").Value = -1280.0588
$ws.Range("N34").Value = -3946.1372
$ws.Range("H99").Value = 2265.6843
$ws.Range("I99").Value = 2153.75
$ws.Range("J99").Value = 2457.5715
$ws.Range("K99").Value = 2153.75
$ws.Range("L99").Value = 2457.5715
$ws.Range("M99").Value = -655.75
$ws.Range("N99").Value = -5453.5715
$ws.Range("H126").Value = 2265.6843
$ws.Range("I126").Value = 2153.75
$ws.Range("J126").Value = 2457.5715
$ws.Range("K126").Value = 6461.25
$ws.Range("L126").Value = 7372.7145
$ws.Range("M126").Value = -3991.25
$ws.Range("N126").Value = -12312.7145
$ws.Range("H132").Value = 2871.85
$ws.Range("I132").Value = 2855.7646
$ws.Range("K132").Value = 8567.293799999999
$ws.Range("M132").Value = -6037.293799999999
$ws.Range("H140").Value = 76999
$ws.Range("J140").Value = 76999
$ws.Range("L140").Value = 76999
$ws.Range("N140").Value = -87359

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 8.666667
$ws.Range("H34").Value = 1741.2106
$ws.Range("J34").Value = 1925.4706
$ws.Range("L34").Value = 5776.4118
$ws.Range("N34").Value = -5944.4118
$ws.Range("H92").Value = 220.6
$ws.Range("I92").Value = 250
$ws.Range("K92").Value = 750
$ws.Range("M92").Value = 498
$ws.Range("H132").Value = 1450
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H17").Value = 600
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -936
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H99").Value = 7119
$ws.Range("I99").Value = 1399
$ws.Range("K99").Value = 1399
$ws.Range("M99").Value = 847
$ws.Range("H102").Value = 2279.7222
$ws.Range("I102").Value = 1698.8334
$ws.Range("K102").Value = 1698.8334
$ws.Range("M102").Value = -76.83339999999998
$ws.Range("H132").Value = 4583.973
$ws.Range("I132").Value = 4163.6
$ws.Range("K132").Value = 12490.8
$ws.Range("M132").Value = -9960.800000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 24419.334
$ws.Range("I99").Value = 24419.334
$ws.Range("K99").Value = 24419.334
$ws.Range("M99").Value = -21424.334
$ws.Range("H100").Value = 2200.3076
$ws.Range("I100").Value = 1850.4
$ws.Range("K100").Value = 1850.4
$ws.Range("M100").Value = -1309.4
$ws.Range("H132").Value = 7400.5884
$ws.Range("I132").Value = 3420
$ws.Range("K132").Value = 10260
$ws.Range("M132").Value = -7730

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7412664.5
$ws.Range("I81").Value = 2387.6667
$ws.Range("J81").Value = 22233218
$ws.Range("K81").Value = 4775.3334
$ws.Range("L81").Value = 44466436
$ws.Range("M81").Value = -3714.3334
$ws.Range("N81").Value = -44468558
$ws.Range("H84").Value = 7412664.5
$ws.Range("I84").Value = 2387.6667
$ws.Range("J84").Value = 22233218
$ws.Range("K84").Value = 23876.667
$ws.Range("L84").Value = 222332180
$ws.Range("M84").Value = -18572.667
$ws.Range("N84").Value = -222342788
$ws.Range("H96").Value = 2494.4443
$ws.Range("J96").Value = 2698.25
$ws.Range("L96").Value = 2698.25
$ws.Range("N96").Value = -5444.25
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
